$d = $word.ActiveDocument

# The document mentions "base_container_job_dn_story_test.yaml" twice:
#   1) In the "Create a job base image" section (the base image build config) - unchanged.
#   2) In the "Create Persistent Volume Claim" section (the PVC creation command) - this
#      one is being renamed to the dedicated PVC creator template file.
#
# Anchor on the PVC section's instruction paragraph so only the second occurrence
# (the one that actually follows it) gets updated.
$anchor = $d.Content
$anchor.Find.Execute("Create PVC using", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$scoped = $d.Range($anchor.End, $d.Content.End)
$scoped.Find.Execute("base_container_job_dn_story_test.yaml", $true, $false, $false, $false, $false, $true, 1, $false, "pvc_creator.yaml", 2)
